# "Literatura i grafovi fix"
# Localises the LOW/MEDIUM/HIGH/EXTREME level labels and the "INCREASE IN %"
# caption to Croatian, and adds axis titles to the line/bar charts on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Translate the workload-level header cells (rows 18, 22, 25) and the
#    "INCREASE IN %" label (row 26) to Croatian. DIFFERENCE (row 23) is left
#    untouched - only the shared-string slot it points at shifts.
# ---------------------------------------------------------------------------
$ws.Range("B18").Value = "Nisko"
$ws.Range("C18").Value = "Srednje"
$ws.Range("D18").Value = "Visoko"
$ws.Range("E18").Value = "Ekstremno"

$ws.Range("B22").Value = "Nisko"
$ws.Range("C22").Value = "Srednje"
$ws.Range("D22").Value = "Visoko"
$ws.Range("E22").Value = "Ekstremno"

$ws.Range("B25").Value = "Nisko"
$ws.Range("C25").Value = "Srednje"
$ws.Range("D25").Value = "Visoko"
$ws.Range("E25").Value = "Ekstremno"

$ws.Range("A26").Value = "Poboljšanje performansa u %"

# ---------------------------------------------------------------------------
# 2. Add axis titles to the six charts on Sheet1.
#    ChartObjects(1..4) are the per-level line charts (Izvođenje / Prosječan
#    broj osvježavanja po sekundi). ChartObjects(5) is the bar chart
#    (Razina opterećenja / Prosječan broj osvježavanja po sekundi).
#    ChartObjects(6) is the "increase %" line chart (Razina opterećenja /
#    Poboljšanje u %).
# ---------------------------------------------------------------------------

function Set-AxisTitles {
    param($chart, [string]$catTitle, [string]$valTitle)

    $catAx = $chart.Axes(1)
    $catAx.HasTitle = $true
    $catAx.AxisTitle.Text = $catTitle

    $valAx = $chart.Axes(2)
    $valAx.HasTitle = $true
    $valAx.AxisTitle.Text = $valTitle
}

# Chart 1 (LOW)
$chart1 = $ws.ChartObjects(1).Chart
Set-AxisTitles $chart1 "Izvođenje" "Prosječan broj osvježavanja po sekundi"

# Chart 2 (MED)
$chart2 = $ws.ChartObjects(2).Chart
Set-AxisTitles $chart2 "Izvođenje" "Prosječan broj osvježavanja po sekundi"

# Chart 3 (HIGH) - note the original author's typo carried over from the
# commit ("ovježavanja" instead of "osvježavanja").
$chart3 = $ws.ChartObjects(3).Chart
Set-AxisTitles $chart3 "Izvođenje" "Prosječan broj ovježavanja po sekundi"

# Chart 4 (EXTREME)
$chart4 = $ws.ChartObjects(4).Chart
Set-AxisTitles $chart4 "Izvođenje" "Prosječan broj osvježavanja po sekundi"

# Chart 5 (bar chart summarising all levels)
$chart5 = $ws.ChartObjects(5).Chart
Set-AxisTitles $chart5 "Razina opterećenja" "Prosječan broj osvježavanja po sekundi"

# Chart 6 (increase-in-% line chart)
$chart6 = $ws.ChartObjects(6).Chart
Set-AxisTitles $chart6 "Razina opterećenja" "Poboljšanje u %"

# ---------------------------------------------------------------------------
# 3. Restore the scroll position / selection that was active when the
#    workbook was saved.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K34").Select()

Write-Output "Literatura i grafovi fix applied"
